# "relative paths to images": the Image column (O) held absolute
# C:\Users\...\images\*.jpg paths baked in by the original author's
# machine. Replace them with paths relative to the workbook location.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$images = [ordered]@{
    2 = "Accord"
    3 = "CRV"
    4 = "Highlander"
    5 = "ModelY"
    6 = "QX60"
    7 = "Lyriq"
    8 = "Pilot"
    9 = "RX350"
}

foreach ($row in $images.Keys) {
    $name = $images[$row]
    $ws.Cells.Item($row, 15).Value = ".\images\$name.jpg"
}

# Match the author's final selection (scrolled right, cell O12 selected).
$ws.Range("O12").Select()
